# Slight touch up with the diagrams:
#   - Slide 1, "Arc 3": give the (previously invisible) head end a visible
#     stealth arrowhead, matching the tail end.
#   - Slide 2, "Arc 20": reposition/resize/rotate the arc and switch its
#     head end from a bare "none" arrow (lg/sm) to a stealth arrow (med/med).
#   - Slide 2, "Arc 3": give the (previously invisible) tail end a visible
#     stealth arrowhead, matching the head end.

$p = $ppt.ActivePresentation

# msoArrowheadStealth = 4
$msoArrowheadStealth = 4
$msoArrowheadWidthMedium = 2
$msoArrowheadLengthMedium = 2

# --- Slide 1 : "Arc 3" -----------------------------------------------------
$s1 = $p.Slides.Item(1)
$arc3Slide1 = $s1.Shapes.Item(24)
$arc3Slide1.Line.BeginArrowheadStyle = $msoArrowheadStealth
# re-touch the (unchanged) tail end so it keeps serializing after the head
# end, i.e. preserve the canonical <a:headEnd/><a:tailEnd/> element order
$arc3Slide1.Line.EndArrowheadStyle = $arc3Slide1.Line.EndArrowheadStyle

# --- Slide 2 : "Arc 20" -----------------------------------------------------
$s2 = $p.Slides.Item(2)
$arc20 = $s2.Shapes.Item(9)
$arc20.Rotation = 120.62495
$arc20.Left = 487.7242913385827
$arc20.Top = 312.61090551181104
$arc20.Width = 24.220984251968503
$arc20.Height = 42.57996062992126
$arc20.Line.BeginArrowheadStyle = $msoArrowheadStealth
$arc20.Line.BeginArrowheadWidth = $msoArrowheadWidthMedium
$arc20.Line.BeginArrowheadLength = $msoArrowheadLengthMedium
# re-touch the (unchanged) tail end so it keeps serializing after the head
# end, i.e. preserve the canonical <a:headEnd/><a:tailEnd/> element order
$arc20.Line.EndArrowheadStyle = $arc20.Line.EndArrowheadStyle

# --- Slide 2 : "Arc 3" -----------------------------------------------------
$arc3Slide2 = $s2.Shapes.Item(24)
$arc3Slide2.Line.EndArrowheadStyle = $msoArrowheadStealth
